$p = $ppt.ActivePresentation

# Slide 2: TextBox "The Moon" -- consolidate "The" / " " / "Moon" runs into one run.
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange
$tr2.Text = "zzzzzzzzzzzzzzzz"
$tr2.Text = "The Moon"

# Slide 3: Title "One More" -- consolidate "One" / " " / "More" runs into one run.
$s3 = $p.Slides.Item(3)
$trTitle = $s3.Shapes.Item(1).TextFrame.TextRange
$trTitle.Text = "zzzzzzzzzzzzzzzz"
$trTitle.Text = "One More"

# Slide 3: TextBox "The Moon" -- consolidate "The" / " " / "Moon" runs into one run.
$trMoon3 = $s3.Shapes.Item(3).TextFrame.TextRange
$trMoon3.Text = "zzzzzzzzzzzzzzzz"
$trMoon3.Text = "The Moon"
